# Add a new "2022-Q1" worksheet (placed right before the "总计" summary sheet)
# and insert its corresponding summary row into "总计".

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")

# --- Create the new "2022-Q1" sheet, positioned before "总计" ---
$new = $wb.Worksheets.Add($total)
$new.Name = "2022-Q1"

# Copy header row formatting (bold font / border / centered) from an existing sheet
$src.Range("B1:H1").Copy()
$new.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Copy the column-A "index" cell formatting down for the 4 data rows
$src.Range("A2").Copy()
$new.Range("A2:A5").PasteSpecial(-4122)  # xlPasteFormats

# The fund code / name / scale / position figures are stored as text (so that
# things like leading zeros and trailing decimal zeros are preserved exactly)
$new.Range("B2:G5").NumberFormat = "@"

# Header row
$new.Cells.Item(1,2).Value = "基金代码"
$new.Cells.Item(1,3).Value = "基金名称"
$new.Cells.Item(1,4).Value = "基金规模"
$new.Cells.Item(1,5).Value = "股票总仓位"
$new.Cells.Item(1,6).Value = "仓位占比"
$new.Cells.Item(1,7).Value = "持有市值(亿元)"
$new.Cells.Item(1,8).Value = "仓位排名"

# Row 2
$new.Cells.Item(2,1).Value = 0
$new.Cells.Item(2,2).Value = "009387"
$new.Cells.Item(2,3).Value = "嘉实稳福混合A"
$new.Cells.Item(2,4).Value = "0.08"
$new.Cells.Item(2,5).Value = "34.71"
$new.Cells.Item(2,6).Value = "3.77"
$new.Cells.Item(2,7).Value = "0.0030"
$new.Cells.Item(2,8).Value = 4

# Row 3
$new.Cells.Item(3,1).Value = 1
$new.Cells.Item(3,2).Value = "009649"
$new.Cells.Item(3,3).Value = "嘉实精选平衡混合A"
$new.Cells.Item(3,4).Value = "0.06"
$new.Cells.Item(3,5).Value = "67.70"
$new.Cells.Item(3,6).Value = "4.89"
$new.Cells.Item(3,7).Value = "0.0029"
$new.Cells.Item(3,8).Value = 4

# Row 4
$new.Cells.Item(4,1).Value = 2
$new.Cells.Item(4,2).Value = "009650"
$new.Cells.Item(4,3).Value = "嘉实精选平衡混合C"
$new.Cells.Item(4,4).Value = "0.01"
$new.Cells.Item(4,5).Value = "67.70"
$new.Cells.Item(4,6).Value = "4.89"
$new.Cells.Item(4,7).Value = "0.0005"
$new.Cells.Item(4,8).Value = 4

# Row 5
$new.Cells.Item(5,1).Value = 3
$new.Cells.Item(5,2).Value = "009388"
$new.Cells.Item(5,3).Value = "嘉实稳福混合C"
$new.Cells.Item(5,4).Value = "0.01"
$new.Cells.Item(5,5).Value = "34.71"
$new.Cells.Item(5,6).Value = "3.77"
$new.Cells.Item(5,7).Value = "0.0004"
$new.Cells.Item(5,8).Value = 4

# --- Update the "总计" sheet: insert a new top row for 2022-Q1, shifting the
#     existing 2021-Q4 / 2021-Q3 rows down by one ---

# Re-fetch the "总计" sheet by name: after inserting "2022-Q1" above, any
# previously-held reference to it now points at the wrong sheet.
$total = $wb.Worksheets.Item("总计")

# Make row 4 match row 3's formatting (so the newly used row picks up the
# same column-A index style) before writing values into it.
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)  # xlPasteFormats

# Shift 2021-Q3 (was row 3) down to row 4
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2021-Q3"
$total.Cells.Item(4,3).Value = 1
$total.Cells.Item(4,4).Value = 0.06

# Shift 2021-Q4 (was row 2) down to row 3
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2021-Q4"
$total.Cells.Item(3,3).Value = 2
$total.Cells.Item(3,4).Value = 1.01

# New row 2 for 2022-Q1
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 4
$total.Cells.Item(2,4).Value = 0.01
